$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Documents")

# --- Update I7: hyperlink now points to the new Payment-Certificate PDF ---
$i7Url = "https://dubaiholding-my.sharepoint.com/personal/arun_naidu_dhre_ae/Documents/Shared/ME/assets/DMS%20149250/PC/PC-04.pdf"
$ws.Range("I7").Value = $i7Url

# --- New row 10: PC#10 payment certificate ---
$ws.Range("A10").Value = "PC10.0"
$ws.Range("B10").Value = "DMS 148857-1"
$ws.Range("C10").Value = "PC#10 - Payment Certificate.pdf"
$ws.Range("D10").Value = "PAYMENT CERTIFICATE"
$ws.Range("E10").Value = 45829
$ws.Range("F10").Value = "PC # 10.0"
$ws.Range("G10").Value = "Payment Certificate # 10"
$ws.Range("H10").Value = "PC-10.PDF"
$i10Url = "https://dubaiholding-my.sharepoint.com/personal/arun_naidu_dhre_ae/Documents/Shared/ME/assets/DMS%20148857-1/PC/PC-09.pdf"
$ws.Range("I10").Value = $i10Url
$ws.Range("I10").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("I10"), $i10Url)
$ws.Range("J10").Value = 45829
$ws.Range("K10").Value = 0

# --- New row 11: PC#11 payment certificate ---
$ws.Range("A11").Value = "PC11.0"
$ws.Range("B11").Value = "DMS 148857-1"
$ws.Range("C11").Value = "PC#11 - Payment Certificate.pdf"
$ws.Range("D11").Value = "PAYMENT CERTIFICATE"
$ws.Range("E11").Value = 45830
$ws.Range("F11").Value = "PC # 11.0"
$ws.Range("G11").Value = "Payment Certificate # 11"
$ws.Range("H11").Value = "PC-11.PDF"
$i11Url = "https://dubaiholding-my.sharepoint.com/personal/arun_naidu_dhre_ae/Documents/Shared/ME/assets/DMS%20148857-1/PC/PC-10.pdf"
$ws.Range("I11").Value = $i11Url
$ws.Range("I11").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("I11"), $i11Url)
$ws.Range("J11").Value = 45830
$ws.Range("K11").Value = 0

# --- Cosmetic: narrower FilePathOrURL column, final selection on the sheet ---
$ws.Columns.Item(9).ColumnWidth = 19
$ws.Range("I11").Select()
